# Rename the "Calendar" sheet to "calendar" (lower-case) as expected by the
# Python script that consumes the exported .csv files. Renaming the sheet
# automatically updates every defined name that refers to it (ActualPCN,
# BestNumberOfEvents, Calendar, CalendarEventNumbers, ... WomPCN) so their
# formula text switches from "Calendar!..." to "calendar!...".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calendar")
$ws.Name = "calendar"

# The sheet-scoped "_xlnm.Print_Area" defined name for this sheet keeps the
# old capitalisation, so refresh it explicitly to match.
$ws.PageSetup.PrintArea = "B1:H30"

# Make the (renamed) calendar sheet the active tab/sheet, and move the
# selection from C298 to F32 (the sheet had been scrolled down to A39 with
# C298 selected; now it opens at the top with F32 selected).
$ws.Activate()
$ws.Range("F32").Select() | Out-Null
